# Update TPM-derived NATMI ligand-receptor metrics (Mdk-Sdc1) with refreshed
# expression/specificity values, per "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value2 = 1.988074333333333
$ws.Range("H2").Value2 = 5.964223
$ws.Range("I2").Value2 = 0.01657769708907969
$ws.Range("J2").Value2 = 0.01657769708907968
$ws.Range("M2").Value2 = 2.565830333333333
$ws.Range("N2").Value2 = 7.697490999999999
$ws.Range("O2").Value2 = 0.0934185609347503
$ws.Range("P2").Value2 = 0.0934185609347503
$ws.Range("Q2").Value2 = 5.10106142938811
$ws.Range("R2").Value2 = 45.90955286449299
$ws.Range("S2").Value2 = 0.001548664605674023
$ws.Range("T2").Value2 = 0.001548664605674023
$ws.Range("G3").Value2 = 1.988074333333333
$ws.Range("H3").Value2 = 5.964223
$ws.Range("I3").Value2 = 0.01657769708907969
$ws.Range("J3").Value2 = 0.01657769708907968
$ws.Range("O3").Value2 = 0.3847798091300315
$ws.Range("P3").Value2 = 0.3847798091300315
$ws.Range("Q3").Value2 = 21.01065809107744
$ws.Range("R3").Value2 = 189.095922819697
$ws.Range("S3").Value2 = 0.00637876312175156
$ws.Range("T3").Value2 = 0.006378763121751558
$ws.Range("G4").Value2 = 1.988074333333333
$ws.Range("H4").Value2 = 5.964223
$ws.Range("I4").Value2 = 0.01657769708907969
$ws.Range("J4").Value2 = 0.01657769708907968
$ws.Range("M4").Value2 = 13.68376133333333
$ws.Range("N4").Value2 = 41.051284
$ws.Range("O4").Value2 = 0.4982080363333638
$ws.Range("P4").Value2 = 0.4982080363333638
$ws.Range("Q4").Value2 = 27.20433469025911
$ws.Range("R4").Value2 = 244.839012212332
$ws.Range("S4").Value2 = 0.008259141913679713
$ws.Range("T4").Value2 = 0.008259141913679711
$ws.Range("G5").Value2 = 1.988074333333333
$ws.Range("H5").Value2 = 5.964223
$ws.Range("I5").Value2 = 0.01657769708907969
$ws.Range("J5").Value2 = 0.01657769708907968
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.6480206666666667
$ws.Range("N5").Value2 = 1.944062
$ws.Range("O5").Value2 = 0.0235935936018545
$ws.Range("P5").Value2 = 0.0235935936018545
$ws.Range("Q5").Value2 = 1.288313254869556
$ws.Range("R5").Value2 = 11.594819293826
$ws.Range("S5").Value2 = 0.0003911274479743924
$ws.Range("T5").Value2 = 0.0003911274479743923
$ws.Range("I6").Value2 = 0.7746030815641455
$ws.Range("J6").Value2 = 0.7746030815641454
$ws.Range("M6").Value2 = 2.565830333333333
$ws.Range("N6").Value2 = 7.697490999999999
$ws.Range("O6").Value2 = 0.0934185609347503
$ws.Range("P6").Value2 = 0.0934185609347503
$ws.Range("Q6").Value2 = 238.3502293002382
$ws.Range("R6").Value2 = 2145.152063702144
$ws.Range("S6").Value2 = 0.07236230517534549
$ws.Range("T6").Value2 = 0.07236230517534548
$ws.Range("I7").Value2 = 0.7746030815641455
$ws.Range("J7").Value2 = 0.7746030815641454
$ws.Range("O7").Value2 = 0.3847798091300315
$ws.Range("P7").Value2 = 0.3847798091300315
$ws.Range("S7").Value2 = 0.2980516258757861
$ws.Range("T7").Value2 = 0.298051625875786
$ws.Range("I8").Value2 = 0.7746030815641455
$ws.Range("J8").Value2 = 0.7746030815641454
$ws.Range("M8").Value2 = 13.68376133333333
$ws.Range("N8").Value2 = 41.051284
$ws.Range("O8").Value2 = 0.4982080363333638
$ws.Range("P8").Value2 = 0.4982080363333638
$ws.Range("Q8").Value2 = 1271.139252318606
$ws.Range("R8").Value2 = 11440.25327086746
$ws.Range("S8").Value2 = 0.3859134802038454
$ws.Range("T8").Value2 = 0.3859134802038453
$ws.Range("I9").Value2 = 0.7746030815641455
$ws.Range("J9").Value2 = 0.7746030815641454
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 0.6480206666666667
$ws.Range("N9").Value2 = 1.944062
$ws.Range("O9").Value2 = 0.0235935936018545
$ws.Range("P9").Value2 = 0.0235935936018545
$ws.Range("Q9").Value2 = 60.19722835322312
$ws.Range("R9").Value2 = 541.7750551790079
$ws.Range("S9").Value2 = 0.0182756703091686
$ws.Range("T9").Value2 = 0.0182756703091686
$ws.Range("G10").Value2 = 23.741365
$ws.Range("H10").Value2 = 71.22409500000001
$ws.Range("I10").Value2 = 0.1979690350870239
$ws.Range("J10").Value2 = 0.1979690350870239
$ws.Range("M10").Value2 = 2.565830333333333
$ws.Range("N10").Value2 = 7.697490999999999
$ws.Range("O10").Value2 = 0.0934185609347503
$ws.Range("P10").Value2 = 0.0934185609347503
$ws.Range("Q10").Value2 = 60.91631447173833
$ws.Range("R10").Value2 = 548.246830245645
$ws.Range("S10").Value2 = 0.01849398236747087
$ws.Range("T10").Value2 = 0.01849398236747087
$ws.Range("G11").Value2 = 23.741365
$ws.Range("H11").Value2 = 71.22409500000001
$ws.Range("I11").Value2 = 0.1979690350870239
$ws.Range("J11").Value2 = 0.1979690350870239
$ws.Range("O11").Value2 = 0.3847798091300315
$ws.Range("P11").Value2 = 0.3847798091300315
$ws.Range("Q11").Value2 = 250.9069677460783
$ws.Range("R11").Value2 = 2258.162709714705
$ws.Range("S11").Value2 = 0.07617448753444157
$ws.Range("T11").Value2 = 0.07617448753444157
$ws.Range("G12").Value2 = 23.741365
$ws.Range("H12").Value2 = 71.22409500000001
$ws.Range("I12").Value2 = 0.1979690350870239
$ws.Range("J12").Value2 = 0.1979690350870239
$ws.Range("M12").Value2 = 13.68376133333333
$ws.Range("N12").Value2 = 41.051284
$ws.Range("O12").Value2 = 0.4982080363333638
$ws.Range("P12").Value2 = 0.4982080363333638
$ws.Range("Q12").Value2 = 324.8711723875534
$ws.Range("R12").Value2 = 2923.840551487981
$ws.Range("S12").Value2 = 0.09862976422551699
$ws.Range("T12").Value2 = 0.09862976422551699
$ws.Range("G13").Value2 = 23.741365
$ws.Range("H13").Value2 = 71.22409500000001
$ws.Range("I13").Value2 = 0.1979690350870239
$ws.Range("J13").Value2 = 0.1979690350870239
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 0.6480206666666667
$ws.Range("N13").Value2 = 1.944062
$ws.Range("O13").Value2 = 0.0235935936018545
$ws.Range("P13").Value2 = 0.0235935936018545
$ws.Range("Q13").Value2 = 15.38489517487667
$ws.Range("R13").Value2 = 138.46405657389
$ws.Range("S13").Value2 = 0.004670800959594516
$ws.Range("T13").Value2 = 0.004670800959594516
$ws.Range("G14").Value2 = 1.301204666666667
$ws.Range("H14").Value2 = 3.903614
$ws.Range("I14").Value2 = 0.01085018625975097
$ws.Range("J14").Value2 = 0.01085018625975097
$ws.Range("M14").Value2 = 2.565830333333333
$ws.Range("N14").Value2 = 7.697490999999999
$ws.Range("O14").Value2 = 0.0934185609347503
$ws.Range("P14").Value2 = 0.0934185609347503
$ws.Range("Q14").Value2 = 3.338670403608222
$ws.Range("R14").Value2 = 30.048033632474
$ws.Range("S14").Value2 = 0.001013608786259937
$ws.Range("T14").Value2 = 0.001013608786259937
$ws.Range("G15").Value2 = 1.301204666666667
$ws.Range("H15").Value2 = 3.903614
$ws.Range("I15").Value2 = 0.01085018625975097
$ws.Range("J15").Value2 = 0.01085018625975097
$ws.Range("O15").Value2 = 0.3847798091300315
$ws.Range("P15").Value2 = 0.3847798091300315
$ws.Range("Q15").Value2 = 13.75158156788289
$ws.Range("R15").Value2 = 123.764234110946
$ws.Range("S15").Value2 = 0.004174932598052269
$ws.Range("T15").Value2 = 0.004174932598052268
$ws.Range("G16").Value2 = 1.301204666666667
$ws.Range("H16").Value2 = 3.903614
$ws.Range("I16").Value2 = 0.01085018625975097
$ws.Range("J16").Value2 = 0.01085018625975097
$ws.Range("M16").Value2 = 13.68376133333333
$ws.Range("N16").Value2 = 41.051284
$ws.Range("O16").Value2 = 0.4982080363333638
$ws.Range("P16").Value2 = 0.4982080363333638
$ws.Range("Q16").Value2 = 17.80537410448623
$ws.Range("R16").Value2 = 160.248366940376
$ws.Range("S16").Value2 = 0.005405649990321777
$ws.Range("T16").Value2 = 0.005405649990321776
$ws.Range("G17").Value2 = 1.301204666666667
$ws.Range("H17").Value2 = 3.903614
$ws.Range("I17").Value2 = 0.01085018625975097
$ws.Range("J17").Value2 = 0.01085018625975097
$ws.Range("K17").Value2 = 3
$ws.Range("L17").Value2 = 1
$ws.Range("M17").Value2 = 0.6480206666666667
$ws.Range("N17").Value2 = 1.944062
$ws.Range("O17").Value2 = 0.0235935936018545
$ws.Range("P17").Value2 = 0.0235935936018545
$ws.Range("Q17").Value2 = 0.8432075155631112
$ws.Range("R17").Value2 = 7.588867640068
$ws.Range("S17").Value2 = 0.0002559948851169901
$ws.Range("T17").Value2 = 0.0002559948851169901
